# Data updated by GitHub Bot (2020-06-06 12:09)
# Applies the 2020-06-06 FOHM covid-19 data refresh to the workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Antal per dag region" - daily case counts per region
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Row 105 (2020-05-29): two corrected values
$ws1.Cells.Item(105, 14).Value = 58   # N105
$ws1.Cells.Item(105, 16).Value = 7    # P105

# Row 123 (2020-06-04): corrected values
$ws1.Cells.Item(123, 2).Value = 1042  # B123
$ws1.Cells.Item(123, 16).Value = 58   # P123
$ws1.Cells.Item(123, 19).Value = 27   # S123
$ws1.Cells.Item(123, 21).Value = 405  # U123

# Row 124 (2020-06-05): full row of figures replaced/updated
$row124 = @{
    2=1016; 3=10; 4=20; 5=1; 6=43; 7=21; 8=13; 9=42; 10=14; 11=19; 12=14;
    13=34; 14=235; 15=3; 16=34; 17=8; 18=5; 19=21; 20=19; 21=373; 22=46; 23=41
}
foreach ($col in $row124.Keys) {
    $ws1.Cells.Item(124, $col).Value = $row124[$col]
}

# Row 125 (NEW, 2020-06-06): copy formatting down from row 124, then fill values
$ws1.Range("A124:W124").Copy()
$ws1.Range("A125").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Cells.Item(125, 1).Value = 43988  # Statistikdatum 2020-06-06
$row125 = @{
    2=114; 3=0; 4=0; 5=0; 6=1; 7=1; 8=11; 9=9; 10=1; 11=0; 12=1;
    13=1; 14=27; 15=0; 16=1; 17=3; 18=0; 19=1; 20=0; 21=27; 22=0; 23=30
}
foreach ($col in $row125.Keys) {
    $ws1.Cells.Item(125, $col).Value = $row125[$col]
}

# ---------------------------------------------------------------------------
# Sheet 2: "Antal avlidna per dag" - daily deaths
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Cells.Item(80, 2).Value = 37
$ws2.Cells.Item(81, 2).Value = 33
$ws2.Cells.Item(82, 2).Value = 28
$ws2.Cells.Item(83, 2).Value = 38
$ws2.Cells.Item(84, 2).Value = 28
$ws2.Cells.Item(85, 2).Value = 19
$ws2.Cells.Item(87, 2).Value = 10

# Row 88 used to hold the "Uppgift saknas" aggregate bucket; it now becomes a
# dated row (2020-06-05), and a new aggregate row 89 is appended using the
# (mis-typed, as published) label "Uppgift saknaa".
$ws2.Range("A87:B87").Copy()
$ws2.Range("A88").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws2.Cells.Item(88, 1).Value = 43987
$ws2.Cells.Item(88, 2).Value = 2

$ws2.Range("A88:B88").Copy()
$ws2.Range("A89").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws2.Cells.Item(89, 1).Value = "Uppgift saknaa"
$ws2.Cells.Item(89, 2).Value = 10

# ---------------------------------------------------------------------------
# Sheet 3: "Antal intensivvårdade per dag" - daily ICU admissions
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Cells.Item(35, 2).Value = 47
$ws3.Cells.Item(55, 2).Value = 33
$ws3.Cells.Item(76, 2).Value = 13
$ws3.Cells.Item(80, 2).Value = 16
$ws3.Cells.Item(85, 2).Value = 19
$ws3.Cells.Item(92, 2).Value = 11
$ws3.Cells.Item(93, 2).Value = 4

# ---------------------------------------------------------------------------
# Sheet 4: "Totalt antal per region" - running totals per region
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$sheet4Data = @{
    2  = @{ B=274;   C=171.67274475097656 }
    3  = @{ B=1297;  C=450.400390625;       D=61 }
    4  = @{ B=89;    C=149.11369323730469 }
    5  = @{ B=1470;  C=511.5142822265625;   E=112 }
    6  = @{ B=887;   C=265.68978881835938;  D=31;  E=57 }
    7  = @{ B=792;   C=605.45831298828125 }
    8  = @{ B=1719;  C=472.77359008789063 }
    9  = @{ B=405;   C=165.0057373046875 }
    10 = @{ B=916;   C=454.6605224609375 }
    11 = @{ B=495;   C=197.92637634277344 }
    12 = @{ B=1957;  C=142.03524780273438; D=97;  E=200 }
    13 = @{ B=14571; C=612.97869873046875; D=822; E=2137 }
    14 = @{ B=1657;  C=556.89990234375 }
    15 = @{ B=2171;  C=565.7874755859375 }
    16 = @{ B=606;   C=214.57859802246094 }
    17 = @{ B=523;   C=192.46621704101563 }
    18 = @{ B=853;   C=347.67083740234375 }
    19 = @{ B=1436;  C=520.58221435546875 }
    20 = @{ B=7589;  C=439.71746826171875; D=349; E=602 }
    21 = @{ B=1942;  C=637.128662109375 }
    22 = @{ B=2238;  C=480.77853393554688; E=198 }
}
$colMap4 = @{ "B"=2; "C"=3; "D"=4; "E"=5 }
foreach ($r in $sheet4Data.Keys) {
    $cells = $sheet4Data[$r]
    foreach ($colLetter in $cells.Keys) {
        $ws4.Cells.Item([int]$r, $colMap4[$colLetter]).Value = $cells[$colLetter]
    }
}

# ---------------------------------------------------------------------------
# Sheet 5: "Totalt antal per kön" - running totals per gender
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

$ws5.Cells.Item(2, 2).Value = 17493
$ws5.Cells.Item(2, 3).Value = 1595
$ws5.Cells.Item(2, 4).Value = 2554

$ws5.Cells.Item(3, 2).Value = 26393
$ws5.Cells.Item(3, 3).Value = 567
$ws5.Cells.Item(3, 4).Value = 2102

# ---------------------------------------------------------------------------
# Sheet 6: "Totalt antal per åldersgrupp" - running totals per age group
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

$sheet6Data = @{
    2  = @{ B=209 }
    3  = @{ B=580 }
    4  = @{ B=4343; C=81 }
    5  = @{ B=5806; C=96 }
    6  = @{ B=6852; C=246 }
    7  = @{ B=8121; C=568; D=137 }
    8  = @{ B=5218; C=649 }
    9  = @{ B=4391; C=421; D=1022 }
    10 = @{ B=5396; C=86;  D=1914 }
    11 = @{ B=2957; D=1184 }
}
$colMap6 = @{ "B"=2; "C"=3; "D"=4 }
foreach ($r in $sheet6Data.Keys) {
    $cells = $sheet6Data[$r]
    foreach ($colLetter in $cells.Keys) {
        $ws6.Cells.Item([int]$r, $colMap6[$colLetter]).Value = $cells[$colLetter]
    }
}

# ---------------------------------------------------------------------------
# Sheet 7: rename "FOHM  5 Jun 2020" -> "FOHM  6 Jun 2020"
# ---------------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item(7)
$ws7.Name = "FOHM  6 Jun 2020"

# ---------------------------------------------------------------------------
# Final view state: sheet 2 ("Antal avlidna per dag") becomes the active tab,
# with the newly appended row selected - matching the published snapshot.
# ---------------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("B89").Select()
